# promo review fall 2021
# Replace the award-list sheet's data with the updated award entries,
# clear out all the now-stale rows, and resize column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New award rows (2-4) -------------------------------------------------
# Column A ("what") for all three new rows first ...
$ws.Range("A2").Value = "Department Academic Excellence Award"
$ws.Range("A3").Value = "Department Industrious Graduate Student Award"
$ws.Range("A4").Value = "Department Graduate Student Teacher of the Year Award"

# ... then column B ("when") ...
$ws.Range("B2").Value = 2015
$ws.Range("B3").Value = 2014
$ws.Range("B4").Value = 2005

# ... then column D ("where") ...
$ws.Range("D2").Value = "Utah State University"
$ws.Range("D3").Value = "Utah State University"
$ws.Range("D4").Value = "Utah State University"

# ... then column C ("with") ...
$ws.Range("C2").Value = "Mathematics and Statistics Department"
$ws.Range("C3").Value = "Mathematics and Statistics Department"
$ws.Range("C4").Value = "Mathematics and Statistics Department"

# No "why" text for the new rows.
$ws.Range("E2").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("E4").ClearContents()

# --- Remove every other stale award row -----------------------------------
$ws.Range("A5:E5").ClearContents()
$ws.Range("A6:E6").ClearContents()
$ws.Range("A7:E15").ClearContents()
$ws.Range("A16:D16").ClearContents()

# --- Column sizing ----------------------------------------------------------
$ws.Columns("A:A").ColumnWidth = 14.6667
